# Allow n/a for judge non-availability
# Fills in the previously-blank row 10 template row with a real example:
#   B10 = judge name (text), C10 = VLJ # (number), D10 = "N/A" (text)
# B10/D10 need an explicit Text number format so "N/A" (and the name) are
# stored/rendered as literal text rather than being reinterpreted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Lamphere, Doris"

$ws.Range("C10").Value = 861

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "N/A"
